$d = $word.ActiveDocument

$replacements = @(
    @{old = "10÷8=1, 2"; new = "65÷4=16, 1"},
    @{old = "63÷6=10, 3"; new = "46÷9=5, 1"},
    @{old = "20÷5=4, 0"; new = "33÷5=6, 3"},
    @{old = "96÷6=16, 0"; new = "73÷3=24, 1"},
    @{old = "35÷4=8, 3"; new = "80÷8=10, 0"},
    @{old = "93÷5=18, 3"; new = "79÷8=9, 7"},
    @{old = "30÷5=6, 0"; new = "35÷9=3, 8"},
    @{old = "80÷2=40, 0"; new = "32÷5=6, 2"},
    @{old = "76÷8=9, 4"; new = "31÷5=6, 1"},
    @{old = "43÷8=5, 3"; new = "21÷5=4, 1"},
    @{old = "89÷9=9, 8"; new = "99÷4=24, 3"},
    @{old = "95÷8=11, 7"; new = "73÷3=24, 1"},
    @{old = "47÷2=23, 1"; new = "17÷4=4, 1"},
    @{old = "38÷9=4, 2"; new = "59÷8=7, 3"},
    @{old = "76÷6=12, 4"; new = "53÷9=5, 8"},
    @{old = "77÷5=15, 2"; new = "56÷9=6, 2"},
    @{old = "26÷4=6, 2"; new = "90÷8=11, 2"},
    @{old = "29÷6=4, 5"; new = "73÷3=24, 1"},
    @{old = "84÷3=28, 0"; new = "36÷2=18, 0"},
    @{old = "97÷3=32, 1"; new = "16÷4=4, 0"},
    @{old = "24÷8=3, 0"; new = "30÷8=3, 6"},
    @{old = "85÷7=12, 1"; new = "33÷6=5, 3"},
    @{old = "85÷3=28, 1"; new = "78÷4=19, 2"},
    @{old = "26÷7=3, 5"; new = "69÷2=34, 1"},
    @{old = "81÷9=9, 0"; new = "76÷4=19, 0"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
